$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.992.29"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "3.786.73"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("D5").Value = "'602.60"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "'163.41"
$ws.Range("E6").Value = "  -2.02%  "
$ws.Range("D7").Value = "3.783.34"
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("E10").Value = "  -2.11%  "
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").Value = "'6.74"
$ws.Range("E12").Value = "  +6.99%  "
$ws.Range("D13").Value = "'0.0000247"
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("D14").Value = "'35.05"
$ws.Range("E14").Value = "  -2.41%  "
$ws.Range("D15").Value = "4.416.80"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Value = "3.781.17"
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("D17").Value = "67.902.68"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "'18.18"
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("E19").Value = "  +1.92%  "
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D21").Value = "'458.28"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("D22").Value = "'9.46"
$ws.Range("E22").Value = "  -4.56%  "
$ws.Range("D23").Value = "'0.690"
$ws.Range("E23").Value = "  -1.63%  "
$ws.Range("D24").Value = "'83.21"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("E25").Value = "  -3.54%  "
$ws.Range("D26").Value = "'11.87"
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'9.90"
$ws.Range("E29").Value = "  -2.09%  "
$ws.Range("D30").Value = "3.934.16"
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("E31").Value = "  -2.75%  "
$ws.Range("E32").Value = "  -7.04%  "
$ws.Range("E33").Value = "  -2.32%  "
$ws.Range("D34").Value = "'28.99"
$ws.Range("E34").Value = "  -2.32%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "'8.91"
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("D38").Value = "'0.147"
$ws.Range("E38").Value = "  +6.53%  "
$ws.Range("D39").Value = "'5.81"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").Value = "'3.22"
$ws.Range("E40").Value = "  -3.78%  "
$ws.Range("D41").Value = "'0.978"
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "'43.67"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("D45").Value = "'47.11"
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("D46").Value = "'152.39"
$ws.Range("E46").Value = "  +2.78%  "
$ws.Range("D47").Value = "'0.294"
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("D48").Value = "'1.38"
$ws.Range("E48").Value = "  -1.71%  "
$ws.Range("D49").Value = "'8.29"
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("D50").Value = "'1.83"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("D51").Value = "'26.44"
$ws.Range("E51").Value = "  -7.36%  "

# Reset style on cells that required a text-forcing apostrophe prefix,
# so no quotePrefix style flag lingers on the saved cell (keeps cell style
# identical to the untouched neighboring cells: no explicit style).
foreach ($addr in @("D4","D5","D6","D12","D13","D14","D18","D21","D22","D23","D24","D26","D29","D34","D36","D38","D39","D40","D41","D42","D44","D45","D46","D47","D48","D49","D50","D51")) {
    $ws.Range($addr).Style = "Normal"
}
